{"js": "// Insert a new bulleted list item right after the paragraph that reads\n// \"Pesquisar livro por categoria para o requisitar;\" so the list becomes:\n//   ... Pesquisar livro por categoria para o requisitar;\n//   ... Deixar um coment\u00e1rio num livro previamente requisitado;   <- NEW\n//   ... Pesquisar por utilizador para o seguir;\n\nconst anchorText = \"Pesquisar livro por categoria para o requisitar;\";\n\nconst results = context.document.body.search(anchorText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Anchor paragraph not found: \" + anchorText);\n}\n\nconst anchorRange = results.items[0];\n\n// Insert a new paragraph right after the found text; it inherits the\n// anchor paragraph's formatting (list style + numbering) and we set its\n// text directly via insertParagraph.\nanchorRange.insertParagraph(\n  \"Deixar um coment\u00e1rio num livro previamente requisitado;\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "# Insert a new bulleted list item right after the paragraph that reads\n# \"Pesquisar livro por categoria para o requisitar;\" so it becomes:\n#   ... Pesquisar livro por categoria para o requisitar;\n#   ... Deixar um coment\u00e1rio num livro previamente requisitado;   <- NEW\n#   ... Pesquisar por utilizador para o seguir;\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph by its text.\n$anchorText = \"Pesquisar livro por categoria para o requisitar;\"\n$paras = $d.Paragraphs\n$targetIndex = -1\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    $p = $paras.Item($i)\n    if ($p.Range.Text.TrimEnd(\"`r`a`n\") -eq $anchorText) {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Anchor paragraph not found: $anchorText\"\n}\n\n$target = $paras.Item($targetIndex)\n\n# InsertParagraphAfter() creates a new paragraph right after the anchor,\n# inheriting its paragraph formatting (style + list numbering).\n$target.Range.InsertParagraphAfter()\n\n# The freshly created paragraph is now the paragraph right after the anchor.\n$newPara = $d.Paragraphs.Item($targetIndex + 1)\n$newPara.Range.Text = \"Deixar um coment\u00e1rio num livro previamente requisitado;\"\n"}
